# Apply the two changes recorded in the commit:
#  1. Slide 5's table switches from the deck's custom "Table_0" style to the
#     built-in table style {6982D59F-941D-4AA9-8CA6-E766E50FB1B4}.
#  2. The presentation's theme colour scheme changes from the "Integral"
#     design's "Red Violet" palette to the stock "Office Theme" / "Office"
#     palette (dk1/lt1/dk2/lt2/accent1-6/hlink/folHlink).

$p = $ppt.ActivePresentation

# --- 1. Table style on slide 5 ---------------------------------------------
$tableShape = $p.Slides.Item(5).Shapes.Item(2)
$tableShape.Table.ApplyStyle("{6982D59F-941D-4AA9-8CA6-E766E50FB1B4}")

# --- 2. Theme colour scheme swap (Integral/Red Violet -> Office Theme/Office)
# PowerPoint's RGB() returns r | (g << 8) | (b << 16); values below are
# precomputed from the target hex triples noted in each comment.
$colors = $p.SlideMaster.Theme.ThemeColorScheme

$colors.Item(1).RGB  = 0          # dk1      000000
$colors.Item(2).RGB  = 16777215   # lt1      FFFFFF
$colors.Item(3).RGB  = 6968388    # dk2      44546A
$colors.Item(4).RGB  = 15132391   # lt2      E7E6E6
$colors.Item(5).RGB  = 13998939   # accent1  5B9BD5
$colors.Item(6).RGB  = 3243501    # accent2  ED7D31
$colors.Item(7).RGB  = 10855845   # accent3  A5A5A5
$colors.Item(8).RGB  = 49407      # accent4  FFC000
$colors.Item(9).RGB  = 12874308   # accent5  4472C4
$colors.Item(10).RGB = 4697456    # accent6  70AD47
$colors.Item(11).RGB = 12673797   # hlink    0563C1
$colors.Item(12).RGB = 7491477    # folHlink 954F72
